# Update example data CON (DTI) connectivity matrix (subject35).
# A handful of cells were asymmetric relative to their mirrored
# (row,col)<->(col,row) counterpart; this resets each such cell to the
# value already present at its transposed position so the matrix is
# symmetric, matching the corrected connectivity data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 2).Value = 0.77805481981011981
$ws.Cells.Item(1, 11).Value = 0.89323638542107808
$ws.Cells.Item(3, 2).Value = 0.62778220197731294
$ws.Cells.Item(4, 2).Value = 0.92782800040264035
$ws.Cells.Item(4, 3).Value = 0.76863616550765279
$ws.Cells.Item(4, 6).Value = 0.70626750353688916
$ws.Cells.Item(5, 4).Value = 0.95901473989895181
$ws.Cells.Item(5, 7).Value = 0.9654759124649922
$ws.Cells.Item(7, 6).Value = 0.76455817872765008
$ws.Cells.Item(7, 9).Value = 0.91505019002772081
$ws.Cells.Item(8, 6).Value = 0.9875368996607139
$ws.Cells.Item(8, 9).Value = 0.89971453718792738
$ws.Cells.Item(9, 10).Value = 0.67708369202151197
$ws.Cells.Item(9, 11).Value = 0.83590685406778897
$ws.Cells.Item(10, 8).Value = 0.96866755761251244
$ws.Cells.Item(10, 11).Value = 0.97584434008922827
$ws.Cells.Item(10, 13).Value = 0.78740814579998153
$ws.Cells.Item(12, 13).Value = 0.99039997935458113
$ws.Cells.Item(13, 14).Value = 0.86649779571684504
$ws.Cells.Item(13, 15).Value = 0.90241236469586927
$ws.Cells.Item(14, 12).Value = 0.87806459229964529
$ws.Cells.Item(14, 16).Value = 0.89552584468129326
$ws.Cells.Item(15, 14).Value = 0.69101019964954113
$ws.Cells.Item(15, 17).Value = 0.69724972779070926
$ws.Cells.Item(16, 15).Value = 0.98344018233042352
$ws.Cells.Item(16, 17).Value = 0.82234028422645666
$ws.Cells.Item(18, 16).Value = 0.72278078638833509
$ws.Cells.Item(18, 17).Value = 0.72572773057756956
$ws.Cells.Item(18, 19).Value = 0.70844217390144082
$ws.Cells.Item(18, 20).Value = 0.71155467067133193
$ws.Cells.Item(19, 20).Value = 0.95219876891614585
$ws.Cells.Item(20, 22).Value = 0.98632412053164575
$ws.Cells.Item(21, 19).Value = 0.58223007222544121
$ws.Cells.Item(21, 20).Value = 0.71702968762942909
$ws.Cells.Item(21, 22).Value = 0.57510478719965041
$ws.Cells.Item(22, 23).Value = 0.98859201276949038
$ws.Cells.Item(22, 24).Value = 0.63229817965141188
$ws.Cells.Item(23, 21).Value = 0.84989197791490867
$ws.Cells.Item(23, 25).Value = 0.8504692436327197
$ws.Cells.Item(24, 23).Value = 0.91982997610897055
$ws.Cells.Item(24, 26).Value = 0.98163245168551261
$ws.Cells.Item(25, 24).Value = 0.89116520284809275
$ws.Cells.Item(25, 26).Value = 0.95647480769653204
$ws.Cells.Item(26, 28).Value = 0.65354092795664553
$ws.Cells.Item(27, 11).Value = 0.91403783777312231
$ws.Cells.Item(27, 26).Value = 0.88919615022700382
$ws.Cells.Item(27, 28).Value = 0.69364500572990129
$ws.Cells.Item(27, 46).Value = 0.75893515306126424
$ws.Cells.Item(28, 30).Value = 0.96577850036186952
$ws.Cells.Item(29, 11).Value = 0.71806470181260074
$ws.Cells.Item(29, 25).Value = 0.83197208131412403
$ws.Cells.Item(29, 28).Value = 0.92889725198776918
$ws.Cells.Item(29, 30).Value = 0.90843291108291102
$ws.Cells.Item(29, 64).Value = 0.97307998931858364
$ws.Cells.Item(30, 50).Value = 0.88901960749403441
$ws.Cells.Item(31, 33).Value = 0.96891081008368085
$ws.Cells.Item(32, 30).Value = 0.89268808157888724
$ws.Cells.Item(32, 33).Value = 0.86046911216189215
$ws.Cells.Item(32, 34).Value = 0.73151173268323788
$ws.Cells.Item(34, 33).Value = 0.93194697094318379
$ws.Cells.Item(34, 36).Value = 0.84784975958818254
$ws.Cells.Item(35, 33).Value = 0.77355370211990282
$ws.Cells.Item(35, 34).Value = 0.8833327864277154
$ws.Cells.Item(35, 37).Value = 0.83567136410184062
$ws.Cells.Item(36, 25).Value = 0.99699308970143385
$ws.Cells.Item(36, 35).Value = 0.8203465417578728
$ws.Cells.Item(36, 37).Value = 0.82596277171365506
$ws.Cells.Item(37, 11).Value = 0.87992077409304681
$ws.Cells.Item(37, 39).Value = 0.92076314236359813
$ws.Cells.Item(38, 13).Value = 0.58116721864548149
$ws.Cells.Item(38, 39).Value = 0.76927133511676171
$ws.Cells.Item(39, 41).Value = 0.7570904893762469
$ws.Cells.Item(40, 31).Value = 0.97203434602506311
$ws.Cells.Item(40, 42).Value = 0.7677965272306998
$ws.Cells.Item(40, 49).Value = 0.69232578034783288
$ws.Cells.Item(42, 41).Value = 0.82848469535643732
$ws.Cells.Item(43, 41).Value = 0.87292902964628882
$ws.Cells.Item(43, 42).Value = 0.91360947350203037
$ws.Cells.Item(43, 44).Value = 0.77827540971135023
$ws.Cells.Item(43, 45).Value = 0.99415013327913537
$ws.Cells.Item(44, 37).Value = 0.78365666048282567
$ws.Cells.Item(44, 42).Value = 0.73791538116566802
$ws.Cells.Item(44, 45).Value = 0.86171051974776947
$ws.Cells.Item(45, 46).Value = 0.86091598275169878
$ws.Cells.Item(45, 47).Value = 0.9382268512612586
$ws.Cells.Item(46, 47).Value = 0.91523650794033951
$ws.Cells.Item(47, 17).Value = 0.76013722176524756
$ws.Cells.Item(47, 49).Value = 0.8566627038839072
$ws.Cells.Item(48, 47).Value = 0.85574241654638694
$ws.Cells.Item(48, 49).Value = 0.80630852570855449
$ws.Cells.Item(49, 51).Value = 0.90046737218892603
$ws.Cells.Item(49, 56).Value = 0.59804844807271251
$ws.Cells.Item(50, 48).Value = 0.83311274318830675
$ws.Cells.Item(50, 51).Value = 0.65669420450364369
$ws.Cells.Item(51, 52).Value = 0.81906130910741592
$ws.Cells.Item(51, 53).Value = 0.95572977663273084
$ws.Cells.Item(52, 3).Value = 0.85730071503369309
$ws.Cells.Item(52, 50).Value = 0.82928453272054936
$ws.Cells.Item(52, 54).Value = 0.71555660839889668
$ws.Cells.Item(53, 52).Value = 0.95042680978544336
$ws.Cells.Item(53, 55).Value = 0.75091953791175792
$ws.Cells.Item(54, 53).Value = 0.95910972989002552
$ws.Cells.Item(54, 55).Value = 0.94575835979103651
$ws.Cells.Item(54, 56).Value = 0.91377121645006709
$ws.Cells.Item(55, 39).Value = 0.67766385757803627
$ws.Cells.Item(55, 43).Value = 0.99779383564302737
$ws.Cells.Item(56, 55).Value = 0.91011667518770833
$ws.Cells.Item(57, 34).Value = 0.90560483983568396
$ws.Cells.Item(58, 25).Value = 0.84359706214303931
$ws.Cells.Item(58, 56).Value = 0.86616396740024371
$ws.Cells.Item(58, 57).Value = 0.96646914190808308
$ws.Cells.Item(59, 49).Value = 0.7778140844347563
$ws.Cells.Item(59, 61).Value = 0.98763877454165527
$ws.Cells.Item(60, 5).Value = 0.76643505490127195
$ws.Cells.Item(60, 58).Value = 0.96709533900470546
$ws.Cells.Item(60, 61).Value = 0.99001506707921116
$ws.Cells.Item(60, 62).Value = 0.99585572246008769
$ws.Cells.Item(61, 1).Value = 0.80796395212457661
$ws.Cells.Item(61, 40).Value = 0.97678825264502189
$ws.Cells.Item(61, 62).Value = 0.89116071354682236
$ws.Cells.Item(62, 64).Value = 0.81907678123833894
$ws.Cells.Item(63, 61).Value = 0.72587077726891291
$ws.Cells.Item(63, 62).Value = 0.80803690839798814
$ws.Cells.Item(64, 63).Value = 0.83989892281583534
$ws.Cells.Item(65, 34).Value = 0.66570409087206595
$ws.Cells.Item(65, 63).Value = 0.64904776337795056
$ws.Cells.Item(66, 7).Value = 0.90671590771677391
$ws.Cells.Item(66, 64).Value = 0.74667772059010273
$ws.Cells.Item(66, 67).Value = 0.99222680081645942
$ws.Cells.Item(66, 68).Value = 0.88302921469677487
$ws.Cells.Item(67, 1).Value = 0.9218846561884938
$ws.Cells.Item(67, 46).Value = 0.7810317553154853
$ws.Cells.Item(67, 65).Value = 0.86304925386906839
$ws.Cells.Item(68, 1).Value = 0.93614880945372381
$ws.Cells.Item(68, 2).Value = 0.89749774934211679
$ws.Cells.Item(68, 58).Value = 0.90624452968736091
